# correctif problème insertion matière dans bd et mise à jour semestre étudiant
# lors de l'inscription

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Mise à jour du semestre des étudiants : les groupes "2-A"/"2-B"/"2-C"
#    (année 2, "S1" précédent) deviennent "1-A"/"1-B"/"1-C" (nouvelle
#    inscription en semestre 1).
for ($r = 3; $r -le 63; $r++) {
    $cell = $ws.Range("D" + $r)
    $val = $cell.Value()
    if ($val -eq "2-B") {
        $cell.Value = "1-B"
    } elseif ($val -eq "2-C") {
        $cell.Value = "1-C"
    } elseif ($val -eq "2-A") {
        $cell.Value = "1-A"
    }
}

# 2) Correctif de l'import : le numéro étudiant (colonne A) portait
#    l'année 2015 au lieu de 2017, et la moyenne (colonne E) a été
#    recalculée suite au correctif d'insertion des notes de matière.
$data = @{
    3 = @(20170926, 16)
    4 = @(20170927, 11)
    5 = @(20170928, 9)
    6 = @(20170929, 17)
    7 = @(20170930, 14)
    8 = @(20170931, 13)
    9 = @(20170932, 15)
    10 = @(20170933, 5)
    11 = @(20170934, 15)
    12 = @(20170935, 11)
    13 = @(20170936, 8)
    14 = @(20170937, 11)
    15 = @(20170938, 20)
    16 = @(20170939, 15)
    17 = @(20170940, 17)
    18 = @(20170941, 17)
    19 = @(20170942, 17)
    20 = @(20170943, 20)
    21 = @(20170944, 7)
    22 = @(20170945, 6)
    23 = @(20170946, 10)
    24 = @(20170947, 11)
    25 = @(20170948, 10)
    26 = @(20170949, 15)
    27 = @(20170950, 17)
    28 = @(20170951, 6)
    29 = @(20170952, 7)
    30 = @(20170953, 11)
    31 = @(20170954, 10)
    32 = @(20170955, 20)
    33 = @(20170956, 20)
    34 = @(20170957, 9)
    35 = @(20170958, 18)
    36 = @(20170959, 7)
    37 = @(20170960, 16)
    38 = @(20170961, 13)
    39 = @(20170962, 14)
    40 = @(20170963, 10)
    41 = @(20170964, 11)
    42 = @(20170965, 16)
    43 = @(20170966, 18)
    44 = @(20170967, 15)
    45 = @(20170968, 11)
    46 = @(20170969, 20)
    47 = @(20170970, 13)
    48 = @(20170971, 5)
    49 = @(20170972, 5)
    50 = @(20170973, 14)
    51 = @(20170974, 11)
    52 = @(20170975, 20)
    53 = @(20170976, 10)
    54 = @(20170977, 9)
    55 = @(20170978, 7)
    56 = @(20170979, 6)
    57 = @(20170980, 17)
    58 = @(20170981, 7)
    59 = @(20170982, 10)
    60 = @(20170983, 10)
    61 = @(20170984, 15)
    62 = @(20170985, 9)
    63 = @(20170986, 15)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Range("A" + $r).Value = $vals[0]
    $ws.Range("E" + $r).Value = $vals[1]
}
